$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.128.42'
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = '2.995.27'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '506.25'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Range("E5").Value = '  -0.58%  '
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '138.24'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.435'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Range("E8").Value = '  +0.23%  '
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '7.48'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Range("E9").Value = '  -1.51%  '
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '0.109'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Range("E10").Value = '  +0.63%  '
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.364'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("D12").Value = '3.524.03'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("E13").Value = '  +1.35%  '
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '26.29'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Range("E14").Value = '  +2.46%  '
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '0.0000162'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Range("E15").Value = '  +5.11%  '
$ws.Range("D16").Value = '57.182.57'
$ws.Range("E16").Value = '  +0.98%  '
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '6.19'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Range("E17").Value = '  +5.74%  '
$ws.Range("D18").Value = '3.000.30'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '12.68'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Range("E19").Value = '  +1.18%  '
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '7.90'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Range("E20").Value = '  +0.33%  '
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '327.85'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Range("E21").Value = '  +0.12%  '
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '0.998'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Range("E22").Value = '  -0.11%  '
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '0.494'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Range("E23").Value = '  +3.23%  '
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '64.30'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Range("E24").Value = '  +2.77%  '
$ws.Range("E25").Value = '  +0.69%  '
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '1.00'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Range("E26").Value = '  -0.23%  '
$ws.Range("D27").Value = '0.0₃0909'
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("E28").Value = '  +2.01%  '
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '7.37'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Range("E29").Value = '  +5.23%  '
$ws.Range("E30").Value = '  +1.54%  '
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '20.50'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("E33").Value = '  +3.20%  '
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '153.73'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Range("E34").Value = '  -1.59%  '
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '5.84'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Range("E35").Value = '  +3.59%  '
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '1.26'
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Range("E36").Value = '  -0.90%  '
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '24.25'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Range("E37").Value = '  +1.70%  '
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.0677'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").Value = '3.035.00'
$ws.Range("E39").Value = '  -0.45%  '
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '37.15'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Range("E40").Value = '  +1.44%  '
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '3.81'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Range("E42").Value = '  +5.41%  '
$ws.Range("D43").Value = '2.289.12'
$ws.Range("E43").Value = '  +0.88%  '
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '0.648'
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("E45").Value = '  -0.58%  '
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '0.975'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Range("E46").Value = '  -2.96%  '
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '5.99'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Range("E47").Value = '  +3.16%  '
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '0.0237'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Range("E48").Value = '  +0.23%  '
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '19.28'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("E50").Value = '  -8.98%  '
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '0.0889'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Range("E51").Value = '  +1.54%  '
